$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to have a 6th column ("ElementName3") that is no longer
# part of the table. Remove it (this also shrinks the used range from
# A1:F13 down to A1:E13 and the ElementName3 header string disappears).
$ws.Range("F1").ClearContents() | Out-Null

# Fill in the module data (rows 2-13, columns A-E) that was added to the
# table.
$ws.Range("A2").Value = "GSEA31"
$ws.Range("B2").Value = "pede. Suspendisse dui."
$ws.Range("C2").Value = "EL Haddad"
$ws.Range("D2").Value = "Nullam feugiat placerat"
$ws.Range("E2").Value = "varius et, euismod"

$ws.Range("A3").Value = "GSEA32"
$ws.Range("B3").Value = "a nunc. In"
$ws.Range("C3").Value = "Badir"
$ws.Range("D3").Value = "sodales nisi magna"
$ws.Range("E3").Value = "elementum sem, vitae"

$ws.Range("A4").Value = "GSEA33"
$ws.Range("B4").Value = "amet metus. Aliquam"
$ws.Range("C4").Value = "Ezzine"
$ws.Range("D4").Value = "Cras vulputate velit"
$ws.Range("E4").Value = "scelerisque neque sed"

$ws.Range("A5").Value = "GSEA34"
$ws.Range("B5").Value = "quam vel sapien"
$ws.Range("C5").Value = "El Alami Hassoun"
$ws.Range("D5").Value = "Nunc mauris elit,"
$ws.Range("E5").Value = "libero et tristique"

$ws.Range("A6").Value = "GSEA35"
$ws.Range("B6").Value = "feugiat nec, diam."
$ws.Range("C6").Value = "Lazaar"
$ws.Range("D6").Value = "pellentesque. Sed dictum."
$ws.Range("E6").Value = "ridiculus mus. Proin"

$ws.Range("A7").Value = "GSEA36"
$ws.Range("B7").Value = "nonummy. Fusce fermentum"
$ws.Range("C7").Value = "El Haddad"
$ws.Range("D7").Value = "neque pellentesque massa"
$ws.Range("E7").Value = "Mauris eu turpis."

$ws.Range("A8").Value = "GSEA41"
$ws.Range("B8").Value = "a, arcu. Sed"
$ws.Range("C8").Value = "EL Haddad"
$ws.Range("D8").Value = "sit amet risus."
$ws.Range("E8").Value = "Nulla facilisi. Sed"

$ws.Range("A9").Value = "GSEA42"
$ws.Range("B9").Value = "Suspendisse eleifend. Cras"
$ws.Range("C9").Value = "El Alami Hassoun"
$ws.Range("D9").Value = "velit dui, semper"
$ws.Range("E9").Value = "ligula elit, pretium"

$ws.Range("A10").Value = "GSEA43"
$ws.Range("B10").Value = "ante. Nunc mauris"
$ws.Range("C10").Value = "Badir"
$ws.Range("D10").Value = "tortor at risus."
$ws.Range("E10").Value = "felis. Donec tempor,"

$ws.Range("A11").Value = "GSEA44"
$ws.Range("B11").Value = "lobortis quam a"
$ws.Range("C11").Value = "Ezzine"
$ws.Range("D11").Value = "euismod est arcu"
$ws.Range("E11").Value = "ligula eu enim."

$ws.Range("A12").Value = "GSEA45"
$ws.Range("B12").Value = "rhoncus. Nullam velit"
$ws.Range("C12").Value = "Ben Achrab"
$ws.Range("D12").Value = "ut dolor dapibus"
$ws.Range("E12").Value = "commodo tincidunt nibh."

$ws.Range("A13").Value = "GSEA46"
$ws.Range("B13").Value = "Donec tincidunt. Donec"
$ws.Range("C13").Value = "EL Haddad"
$ws.Range("D13").Value = "ornare tortor at"
$ws.Range("E13").Value = "ac, feugiat non,"

# B2 carries an explicit (non-theme) black font color, as if it had been
# pasted in from another source.
$ws.Range("B2").Font.Color = 0

# Resize the columns to fit the new, wider content.
$ws.Columns.Item(1).ColumnWidth = 6.7
$ws.Columns.Item(2).ColumnWidth = 26.2
$ws.Columns.Item(3).ColumnWidth = 15.35
$ws.Columns.Item(4).ColumnWidth = 24.35
$ws.Columns.Item(5).ColumnWidth = 22.7

# Update the current selection.
$ws.Range("H9").Select() | Out-Null
